# "Download file bug fixed" - add a LoginRetry setting so the FannieMae
# login-retry count can be configured alongside the other Settings/Assets.

$wb = $excel.ActiveWorkbook

$wsSettings = $wb.Worksheets.Item("Settings")
$wsAssets   = $wb.Worksheets.Item("Assets")

# Populate the new cells. Assets is filled in first so that the shared
# strings are interned in the same order the workbook expects
# (LoginRetry, FannieMae_LoginRetry, SetLoginRetry).
$wsAssets.Range("A7").Value = "LoginRetry"
$wsAssets.Range("B7").Value = "FannieMae_LoginRetry"

$wsSettings.Range("A4").Value = "SetLoginRetry"
$wsSettings.Range("B4").Value = "FannieMae_LoginRetry"

# Move the selection on each sheet to the newly added row, and make sure
# "Assets" ends up as the active tab again (it was the active sheet in the
# original workbook).
$wsSettings.Activate()
$wsSettings.Range("A4").Select()

$wsAssets.Activate()
$wsAssets.Range("A7").Select()
